# Weekly price-sheet update: insert the latest week's record at row 248
# (pushing the existing rows 248-269 down to 249-270) and populate the
# new row with this week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 248; Excel shifts rows 248:269 down
# to 249:270 and extends the sheet dimension automatically.
$ws.Rows.Item(248).Insert()

# Fill the newly inserted row 248 with the new week's data.
$ws.Cells.Item(248, 1).Value = 8
$ws.Cells.Item(248, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(248, 3).Value = "Coquimbo"
$ws.Cells.Item(248, 4).Value = 45265
$ws.Cells.Item(248, 5).Value = 4
$ws.Cells.Item(248, 6).Value = 100112044
$ws.Cells.Item(248, 7).Value = "Perejil"
$ws.Cells.Item(248, 8).Value = "Sin especificar"
$ws.Cells.Item(248, 9).Value = "Primera"
$ws.Cells.Item(248, 10).Value = 2000
$ws.Cells.Item(248, 11).Value = 2400
$ws.Cells.Item(248, 12).Value = 2500
$ws.Cells.Item(248, 13).Value = 2450
$ws.Cells.Item(248, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(248, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(248, 16).Value = 1633
$ws.Cells.Item(248, 17).Value = 1.5
$ws.Cells.Item(248, 18).Value = "Hortaliza"
